# C5-PowerPoint.pptx edit
#
# 1) The table on slide 6 switches from the custom "Table_0" table
#    style to PowerPoint's built-in "Medium Style 2 - Accent 1" style.
# 2) The deck's theme colour scheme (driving ppt/theme/theme1.xml,
#    which is the Slide Master's theme) is swapped from the "Integral"
#    template palette over to the stock "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
    }
}
$table = $tableShape.Table
$table.ApplyStyle("{B9D60FBA-B6D0-4019-8651-35028C42CB23}")

# --- 2. Theme colours -------------------------------------------------
# Office theme colour scheme values (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink), expressed as VBA-style RGB() long values
# (val = R + G*256 + B*65536).
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

$colorScheme.Colors(1).RGB  = 0          # dk1      000000
$colorScheme.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$colorScheme.Colors(3).RGB  = 6968388    # dk2      44546A
$colorScheme.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$colorScheme.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$colorScheme.Colors(6).RGB  = 3243501    # accent2  ED7D31
$colorScheme.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$colorScheme.Colors(8).RGB  = 49407      # accent4  FFC000
$colorScheme.Colors(9).RGB  = 12874308   # accent5  4472C4
$colorScheme.Colors(10).RGB = 4697456    # accent6  70AD47
$colorScheme.Colors(11).RGB = 12673797   # hlink    0563C1
$colorScheme.Colors(12).RGB = 7491477    # folHlink 954F72
